# Working for Cards filter
# Fix casing of the "MaxFeeFilter" key to "maxFeeFilter" on the main sheet,
# and move the Excel selection to C11 (matching the saved cursor position).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("main")

$ws.Range("B8").Value = "maxFeeFilter"

$ws.Range("C11").Select()
